# msz - field hint and error checks part 1
#
# Adds two new test-data rows (10 & 11) to the "Tabelle1" sheet describing
# hint checks for the vehicle-data page, moves the selection to A9, and
# repositions the screenshot picture that sits below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 10: "Vehicle Page check for hints regarding mandatory fields"
$ws.Cells.Item(10, 1).Value  = "Vehicle Page check for hints regarding mandatory fields"
$ws.Cells.Item(10, 3).Value  = "<HINT Select an option>"
$ws.Cells.Item(10, 4).Value  = "<HINT This field is mandatory>"
$ws.Cells.Item(10, 5).Value  = "<HINT This field is mandatory>"
$ws.Cells.Item(10, 6).Value  = "<HINT Select an option>"
$ws.Cells.Item(10, 7).Value  = "<HINT Select an option>"
$ws.Cells.Item(10, 8).Value  = "<HINT This field is mandatory>"
$ws.Cells.Item(10, 9).Value  = "<NOHINT>"
$ws.Cells.Item(10, 10).Value = "<HINT This field is mandatory>"

# --- New row 11: "Vehicle Page check error hint list price too low"
$ws.Cells.Item(11, 1).Value = "Vehicle Page check error hint list price too low"
$ws.Cells.Item(11, 8).Value = "<HINT Must be a number between 500 and 100000>"

# --- Selection moves to A9
$ws.Range("A9").Select() | Out-Null

# --- Column widths are best-fit (auto-fit) - recompute now that the new,
#     wider strings (columns D/E/H/J) have been entered.
$ws.Columns.AutoFit() | Out-Null

# --- Move the screenshot picture down (it now starts two rows lower,
#     beneath the two newly inserted rows) while keeping its size.
$shp = $ws.Shapes.Item(1)
$shp.Left = 6.6
$shp.Top = 186.6

Write-Host "Applied vehicle-data hint/error rows, selection and picture move"
